$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Rename the worksheet tab: "GST Audit Report" -> "GST Report"
# ---------------------------------------------------------------------------
$ws.Name = "GST Report"

# ---------------------------------------------------------------------------
# 2. Remove the frozen header pane (row 1 is no longer frozen)
# ---------------------------------------------------------------------------
$excel.ActiveWindow.FreezePanes = $false

# ---------------------------------------------------------------------------
# 3. Update the data row (row 2) values before we touch row layout
#    (values that look numeric must stay stored as *text*, exactly like the
#    rest of the sheet, so we force a text value with a leading apostrophe)
# ---------------------------------------------------------------------------
$ws.Range("B2").Value = "WESTSIDE, A UNIT OF TRENT LTD"
$ws.Range("D2").Value = "W089100169940"
$ws.Range("E2").Value = "2024-09-28 17:41:22"
$ws.Range("G2").Value = "'388.06"
$ws.Range("H2").Value = "'194.03"
$ws.Range("I2").Value = "'194.03"
$ws.Range("K2").Value = "996211, 300980061004, 300988526002, 300992658003, 600000562, 300922355001, 300989351001"

# ---------------------------------------------------------------------------
# 4. Remove the TOTAL row (3), blank spacer row (4->5) and footer rows (5,6)
#    Delete from the bottom up so row indices above stay valid.
# ---------------------------------------------------------------------------
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(3).Delete()

# ---------------------------------------------------------------------------
# 5. Strip all the custom cell formatting (fonts/fills/borders/alignment)
#    back to the workbook's plain default "Normal" style.
# ---------------------------------------------------------------------------
$ws.Range("A1:K2").Style = "Normal"

# ---------------------------------------------------------------------------
# 6. Drop the custom row heights so rows 1-2 go back to the sheet default.
# ---------------------------------------------------------------------------
$ws.Rows.Item(1).AutoFit()
$ws.Rows.Item(2).AutoFit()

Write-Host "Applied GST report formatting + content update"
